# Auto-generated script to apply price/quantity updates across crypto sheets
# (values sourced from an external price feed refresh; only cell values change, no formulas/styles)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ETH")
$ws.Range("J3").Value = 2805.530932143699
$ws.Range("B12").Value = 0.00727405
$ws.Range("B35").Value = 0.12557635
$ws.Range("D35").Value = 218.04
$ws.Range("B36").Value = 0.02524552
$ws.Range("D36").Value = 44.5
$ws.Range("B40").Value = 0.05722201
$ws.Range("D40").Value = 107.85

$ws = $wb.Worksheets.Item("APE")
$ws.Range("J3").Value = 1.695708866191765
$ws.Range("B5").Value = 16.34616507
$ws.Range("D5").Value = 44.5
$ws.Range("B6").Value = 0.60016031

$ws = $wb.Worksheets.Item("ATOM")
$ws.Range("J3").Value = 10.38859014094903
$ws.Range("B7").Value = 0.02916115

$ws = $wb.Worksheets.Item("AVAX")
$ws.Range("J3").Value = 40.22166036858602
$ws.Range("B5").Value = 2.65956882
$ws.Range("D5").Value = 44.5
$ws.Range("B6").Value = 0.01672187

$ws = $wb.Worksheets.Item("AMP")
$ws.Range("J3").Value = 0.003831218778778552

$ws = $wb.Worksheets.Item("BNB")
$ws.Range("J3").Value = 353.7808698822783
$ws.Range("B10").Value = 0.00273304
$ws.Range("B11").Value = 0.5819841
$ws.Range("D11").Value = 165.87
$ws.Range("B12").Value = 0.15496513
$ws.Range("D12").Value = 44.5

$ws = $wb.Worksheets.Item("DOGE")
$ws.Range("J3").Value = 0.08400543458557617
$ws.Range("B6").Value = 0.29148391

$ws = $wb.Worksheets.Item("DOT")
$ws.Range("J3").Value = 7.791612979100291
$ws.Range("B5").Value = 7.76131808
$ws.Range("D5").Value = 44.5
$ws.Range("B6").Value = 0.07971715

$ws = $wb.Worksheets.Item("EGLD")
$ws.Range("J3").Value = 58.65352889944671
$ws.Range("B6").Value = 0.00299745

$ws = $wb.Worksheets.Item("GRT")
$ws.Range("J3").Value = 0.2182830066722943

$ws = $wb.Worksheets.Item("ICP")
$ws.Range("J3").Value = 13.60076500324808
$ws.Range("B6").Value = 0.00236028

$ws = $wb.Worksheets.Item("BTC")
$ws.Range("J3").Value = 51816.34204510271
$ws.Range("B6").Value = 0.00035582
$ws.Range("B23").Value = 0.00751975
$ws.Range("D23").Value = 194.25
$ws.Range("B24").Value = 0.00165682
$ws.Range("D24").Value = 44.5
$ws.Range("B34").Value = 0.0020927
$ws.Range("D34").Value = 62.15

$ws = $wb.Worksheets.Item("KAVA")
$ws.Range("J3").Value = 0.7620200559055329

$ws = $wb.Worksheets.Item("LDO")
$ws.Range("J3").Value = 3.177085853972541
$ws.Range("B6").Value = 0.02023803

$ws = $wb.Worksheets.Item("LINK")
$ws.Range("J3").Value = 19.95081045089485
$ws.Range("B6").Value = 0.00248328

$ws = $wb.Worksheets.Item("LTC")
$ws.Range("J3").Value = 70.3674433490243
$ws.Range("B6").Value = 0.00133944

$ws = $wb.Worksheets.Item("LUNA")
$ws.Range("J3").Value = 0.7140387166031454
$ws.Range("B6").Value = 0.05843168

$ws = $wb.Worksheets.Item("LUNC")
$ws.Range("J3").Value = 0.0001285709315128851
$ws.Range("B18").Value = 5027.50637501

$ws = $wb.Worksheets.Item("MATIC")
$ws.Range("J3").Value = 0.9481772358213713
$ws.Range("B6").Value = 0.32916686
$ws.Range("B7").Value = 48.95971409
$ws.Range("D7").Value = 44.5

$ws = $wb.Worksheets.Item("MEME")
$ws.Range("J3").Value = 0.0277240980378696
$ws.Range("B6").Value = 0.06788966

$ws = $wb.Worksheets.Item("MINA")
$ws.Range("J3").Value = 1.382432886133116
$ws.Range("B6").Value = 0.35011212

$ws = $wb.Worksheets.Item("NEAR")
$ws.Range("J3").Value = 3.495103029467528
$ws.Range("B6").Value = 24.0184439
$ws.Range("D6").Value = 44.5
$ws.Range("B7").Value = 0.10300341

$ws = $wb.Worksheets.Item("SEI")
$ws.Range("J3").Value = 0.9652156483141192
$ws.Range("B6").Value = 0.07637707000000001

$ws = $wb.Worksheets.Item("SHIB")
$ws.Range("J3").Value = 0.000009757963151715204
$ws.Range("B6").Value = 276.18

$ws = $wb.Worksheets.Item("SHPING")
$ws.Range("J3").Value = 0.005158418960523025

$ws = $wb.Worksheets.Item("SOL")
$ws.Range("J3").Value = 112.3426479836591
$ws.Range("B16").Value = 6.11236541
$ws.Range("D16").Value = 129.24
$ws.Range("B17").Value = 0.06474977
$ws.Range("B18").Value = 1.91890129
$ws.Range("D18").Value = 44.5

$ws = $wb.Worksheets.Item("TRX")
$ws.Range("J3").Value = 0.1351268544878116
$ws.Range("B6").Value = 0.26546308

$ws = $wb.Worksheets.Item("UNI")
$ws.Range("J3").Value = 7.621863595791873
$ws.Range("B6").Value = 0.00274924

$ws = $wb.Worksheets.Item("XRP")
$ws.Range("J3").Value = 0.5579099297939667
$ws.Range("B6").Value = 0.86605969

$ws = $wb.Worksheets.Item("TIA")
$ws.Range("J3").Value = 18.95302405556108
$ws.Range("B6").Value = 0.0041136

$ws = $wb.Worksheets.Item("DYDX")
$ws.Range("J3").Value = 3.159238440386449
$ws.Range("B6").Value = 0.0008214

$ws = $wb.Worksheets.Item("POLIS")
$ws.Range("J3").Value = 0.4447268980192154

$ws = $wb.Worksheets.Item("ATLAS")
$ws.Range("J3").Value = 0.004871902199245946

$ws = $wb.Worksheets.Item("ACE")
$ws.Range("J3").Value = 10.27029566359748
$ws.Range("B6").Value = 0.00002519

$ws = $wb.Worksheets.Item("ADA")
$ws.Range("J3").Value = 0.6263248237300927
$ws.Range("B6").Value = 0.7821769200000001
$ws.Range("B7").Value = 123.14401613
$ws.Range("D7").Value = 44.5

$ws = $wb.Worksheets.Item("ALGO")
$ws.Range("J3").Value = 0.1950698602832271
$ws.Range("B6").Value = 0.5844629
